$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.265.17'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '2.429.66'
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("E4").Value = '  -0.08%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '559.25'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +0.88%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '161.26'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +0.77%  '
$ws.Range("E7").Value = '  -0.09%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.513'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +2.94%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.156'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +6.89%  '
$ws.Range("E10").Value = '  +0.11%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.326'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -1.93%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '4.77'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +0.71%  '
$ws.Range("D13").Value = '68.302.95'
$ws.Range("E13").Value = '  +0.49%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '0.0000169'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +2.12%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '23.02'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -0.25%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '10.33'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -3.32%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '336.17'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -0.80%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '6.81'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -2.34%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '3.78'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +0.87%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '1.88'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +2.13%  '
$ws.Range("E21").Value = '  +0.08%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '66.77'
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '3.65'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +0.22%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '8.06'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +0.72%  '
$ws.Range("D25").Value = '0.0₃0807'
$ws.Range("E25").Value = '  +0.36%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '7.13'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +1.12%  '
$ws.Range("E27").Value = '  +0.00%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '423.80'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -0.57%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '1.13'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +1.10%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '1.60'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -0.07%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '160.45'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +3.04%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '18.97'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -0.11%  '
$ws.Range("E33").Value = '  -0.13%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '17.73'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +0.34%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.103'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -4.81%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.293'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -2.43%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '4.29'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -0.97%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '1.46'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +0.78%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '1.04'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -3.56%  '
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '3.33'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +1.05%  '
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '1.99'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -0.40%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '128.91'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -1.94%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.0714'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +0.43%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.475'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -0.04%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.556'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -0.07%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.0917'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +1.38%  '
$ws.Range("E47").Value = '  +0.94%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '1.35'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -3.28%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '16.47'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -0.89%  '
$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '4.78'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -6.05%  '
$ws.Range("D51").Value = '0.0₆0202'
$ws.Range("E51").Value = '  +3.16%  '
